# Add git-diff-plugin info & python pip entry to the git cheatsheet (Sheet1).
#
# Target: new rows 29 & 30 filled in with content describing how to wire up
# a custom "git diff" textconv plugin for .xlsx files, plus supporting
# sharedStrings / styles / row-height / selection changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Cell text for the two new rows -----------------------------------
# New shared strings must be created in this exact order so they line up
# with the indices used by the cells that reference them:
#   A29, A30, B30, B29, C29, C30  ->  71, 72, 73, 74, 75, 76

$A29 = 'Config系統'
$A30 = 'Diff'
$B30 = 'Extend git diff for excel'
# B29 is just "/" but Excel stores it with a quote-prefix style (the text
# looks like it could start a path/formula), so type it the way a user
# would in the UI: a leading apostrophe forces text + quotePrefix.
$B29 = "'/"

$C29 = @'
1. File ".gitignore" in a folder is equivalent to "/.git/info/exclude", the only difference is that former is only effective for local folder and latter is valid globally
2. File ".gitattributes" and "/.git/info/attributes" and they are simular relationship like ignore files
'@

$C30 = @'
1. Add entry to  "/.git/info/attributes": 
       $ *.xlsx diff=git_diff_xlsx
2. Add entry to ".git/config":
       $ [diff "git_diff_xlsx"]
       $       binary = True
       $       textconv = python C:/Users/BI77/Documents/playground/git_diff_xlsx.py
3. Then when run git diff *.xlsx, git will use given command to generate difference.
** Interesting fact noticed: git_diff_xlsx.py only need 1 file as input, and output stdin, so that I assume that git has its framework to compare 2 file stream. So by understanding this, I can develope my own plugin too!
** The experiment git repo is the "Playground" folder in my BCM working directory 
'@

# --- Write the values in the order that reproduces the shared-string index order
$ws.Range("A29").Value = $A29
$ws.Range("A30").Value = $A30
$ws.Range("B30").Value = $B30
$ws.Range("B29").Value = $B29
$ws.Range("C29").Value = $C29
$ws.Range("C30").Value = $C30

# --- Row heights (wrap-text rows sized to fit the new content) --------
$ws.Rows.Item(29).RowHeight = 39
$ws.Rows.Item(30).RowHeight = 141

# --- Update the view to where the author ended up scrolled/selected ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 25
$win.ScrollColumn = 1
[void]$ws.Range("C31").Select()
